$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3 (shifts existing rows 3-8 down to 4-9,
# preserving their contents and formatting).
$ws.Rows("3:3").Insert()

# Copy the date cell style (custom date number format) from row 4 (the row
# that was previously row 3) onto the newly inserted row 3's date cell.
$ws.Range("D4").Copy()
$ws.Range("D3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new weekly price record in row 3.
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C3").Value = "Arica y Parinacota"
$ws.Range("D3").Value = 45274
$ws.Range("E3").Value = 15
$ws.Range("F3").Value = 100112017
$ws.Range("G3").Value = "Corazón de apio"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 200
$ws.Range("K3").Value = 800
$ws.Range("L3").Value = 1000
$ws.Range("M3").Value = 900
$ws.Range("N3").Value = "`$/docena de matas"
$ws.Range("O3").Value = "Región de Arica y Parinacota"
$ws.Range("P3").Value = 150
$ws.Range("Q3").Value = 6
$ws.Range("R3").Value = "Hortaliza"

$wb.Save()
